$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Operations")
$ws.Rows("1:10").Insert()
